$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 750.7143
$ws.Range("I33").Value = 789.6923
$ws.Range("J33").Value = 244
$ws.Range("K33").Value = 789.6923
$ws.Range("L33").Value = 244
$ws.Range("M33").Value = -560.6923
$ws.Range("N33").Value = -702
$ws.Range("H41").Value = 2324.75
$ws.Range("I41").Value = 2149.5
$ws.Range("K41").Value = 2149.5
$ws.Range("M41").Value = -1709.5
$ws.Range("H64").Value = 4583.778
$ws.Range("J64").Value = 5222.25
$ws.Range("L64").Value = 5222.25
$ws.Range("N64").Value = -5718.25
$ws.Range("H67").Value = 4583.778
$ws.Range("J67").Value = 5222.25
$ws.Range("L67").Value = 5222.25
$ws.Range("N67").Value = -6938.25
$ws.Range("H101").Value = 25000448
$ws.Range("I101").Value = 25000448
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 75001344
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -74999722
$ws.Range("N101").ClearContents()
$ws.Range("H129").Value = 4036.889
$ws.Range("I129").Value = 4370.6665
$ws.Range("J129").Value = 3369.3333
$ws.Range("K129").Value = 13111.9995
$ws.Range("L129").Value = 10107.9999
$ws.Range("M129").Value = -8111.999500000002
$ws.Range("N129").Value = -20107.9999
$ws.Range("H138").Value = 1925.2162
$ws.Range("I138").Value = 1715.04
$ws.Range("J138").Value = 2363.0833
$ws.Range("K138").Value = 5145.12
$ws.Range("L138").Value = 7089.249899999999
$ws.Range("M138").Value = -5.119999999999891
$ws.Range("N138").Value = -17369.2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 43247.25
$ws.Range("J24").Value = 43247.25
$ws.Range("L24").Value = 43247.25
$ws.Range("N24").Value = -43995.25
$ws.Range("H32").Value = 5175.3184
$ws.Range("I32").Value = 4647.7896
$ws.Range("K32").Value = 4647.7896
$ws.Range("M32").Value = -4360.7896
$ws.Range("H43").Value = 41665
$ws.Range("J43").Value = 39997
$ws.Range("L43").Value = 39997
$ws.Range("N43").Value = -40623
$ws.Range("H61").Value = 5785.4
$ws.Range("I61").Value = 4744.364
$ws.Range("K61").Value = 4744.364
$ws.Range("M61").Value = -4532.364
$ws.Range("H100").Value = 43247.25
$ws.Range("J100").Value = 43247.25
$ws.Range("L100").Value = 43247.25
$ws.Range("N100").Value = -45411.25
$ws.Range("H101").Value = 127598.5
$ws.Range("J101").Value = 127598.5
$ws.Range("L101").Value = 127598.5
$ws.Range("N101").Value = -134088.5
$ws.Range("H122").Value = 1615.1
$ws.Range("I122").Value = 678.8570999999999
$ws.Range("K122").Value = 2036.5713
$ws.Range("M122").Value = 413.4287000000002
$ws.Range("H132").Value = 1048.5714
$ws.Range("I132").Value = 1053
$ws.Range("K132").Value = 3159
$ws.Range("M132").Value = -629
$ws.Range("H136").Value = 5785.4
$ws.Range("I136").Value = 4744.364
$ws.Range("K136").Value = 14233.092
$ws.Range("M136").Value = -11683.092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1583.5
$ws.Range("I20").Value = 1445
$ws.Range("K20").Value = 1445
$ws.Range("M20").Value = -1198
$ws.Range("H134").Value = 1999.8572
$ws.Range("J134").Value = 1950
$ws.Range("L134").Value = 5850
$ws.Range("N134").Value = -10920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3410.2
$ws.Range("I16").Value = 2634
$ws.Range("J16").Value = 4574.5
$ws.Range("K16").Value = 2634
$ws.Range("L16").Value = 4574.5
$ws.Range("M16").Value = -2347
$ws.Range("N16").Value = -5148.5
$ws.Range("H58").Value = 1319.3334
$ws.Range("I58").Value = 1361.3478
$ws.Range("K58").Value = 1361.3478
$ws.Range("M58").Value = -1158.3478
$ws.Range("H99").Value = 4020.8
$ws.Range("I99").Value = 2105
$ws.Range("J99").Value = 4499.75
$ws.Range("K99").Value = 2105
$ws.Range("L99").Value = 4499.75
$ws.Range("M99").Value = -607
$ws.Range("N99").Value = -7495.75
$ws.Range("H113").Value = 3410.2
$ws.Range("I113").Value = 2634
$ws.Range("J113").Value = 4574.5
$ws.Range("K113").Value = 2634
$ws.Range("L113").Value = 4574.5
$ws.Range("M113").Value = -464
$ws.Range("N113").Value = -8914.5
$ws.Range("H126").Value = 4020.8
$ws.Range("I126").Value = 2105
$ws.Range("J126").Value = 4499.75
$ws.Range("K126").Value = 6315
$ws.Range("L126").Value = 13499.25
$ws.Range("M126").Value = -3845
$ws.Range("N126").Value = -18439.25
$ws.Range("H132").Value = 2239.3333
$ws.Range("I132").Value = 2125.5
$ws.Range("K132").Value = 6376.5
$ws.Range("M132").Value = -3846.5
$ws.Range("H136").Value = 1319.3334
$ws.Range("I136").Value = 1361.3478
$ws.Range("K136").Value = 4084.0434
$ws.Range("M136").Value = -1534.0434
$ws.Range("H141").Value = 162331.78
$ws.Range("J141").Value = 176998.88
$ws.Range("L141").Value = 176998.88
$ws.Range("N141").Value = -187358.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 59956
$ws.Range("J37").Value = 59956
$ws.Range("L37").Value = 179868
$ws.Range("N37").Value = -180092
$ws.Range("H116").Value = 1642.6666
$ws.Range("I116").Value = 1642.6666
$ws.Range("K116").Value = 4927.9998
$ws.Range("M116").Value = -1485.9998
$ws.Range("H117").Value = 11392.8
$ws.Range("I117").Value = 387.6
$ws.Range("J117").Value = 22398
$ws.Range("K117").Value = 1162.8
$ws.Range("L117").Value = 67194
$ws.Range("M117").Value = 2279.2
$ws.Range("N117").Value = -74078
$ws.Range("H131").Value = 272786.22
$ws.Range("J131").Value = 360013
$ws.Range("L131").Value = 1080039
$ws.Range("N131").Value = -1090119

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 6250
$ws.Range("I40").Value = 6250
$ws.Range("K40").Value = 6250
$ws.Range("M40").Value = -6099
$ws.Range("H102").Value = 1639
$ws.Range("I102").Value = 1639
$ws.Range("K102").Value = 1639
$ws.Range("M102").Value = -17
$ws.Range("H122").Value = 1922
$ws.Range("I122").Value = 1796
$ws.Range("J122").Value = 2048
$ws.Range("K122").Value = 5388
$ws.Range("L122").Value = 6144
$ws.Range("M122").Value = -2938
$ws.Range("N122").Value = -11044
$ws.Range("H132").Value = 2088.0625
$ws.Range("I132").Value = 1970.1538
$ws.Range("K132").Value = 5910.4614
$ws.Range("M132").Value = -3380.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2936.8
$ws.Range("I136").Value = 2213.4546
$ws.Range("J136").Value = 4926
$ws.Range("K136").Value = 6640.3638
$ws.Range("L136").Value = 14778
$ws.Range("M136").Value = -4090.3638
$ws.Range("N136").Value = -19878

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 19950
$ws.Range("J28").Value = 20000
$ws.Range("L28").Value = 20000
$ws.Range("N28").Value = -20696
$ws.Range("H41").Value = 24831.334
$ws.Range("I41").Value = 24497
$ws.Range("K41").Value = 24497
$ws.Range("M41").Value = -24107
$ws.Range("H75").Value = 90118
$ws.Range("I75").Value = 90118
$ws.Range("K75").Value = 90118
$ws.Range("M75").Value = -89182
$ws.Range("H78").Value = 90118
$ws.Range("I78").Value = 90118
$ws.Range("K78").Value = 270354
$ws.Range("M78").Value = -265674
$ws.Range("H96").Value = 803
$ws.Range("I96").Value = 803
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 803
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 570
$ws.Range("N96").ClearContents()
$ws.Range("H122").Value = 1355
$ws.Range("I122").Value = 1355
$ws.Range("K122").Value = 4065
$ws.Range("M122").Value = -1615
$ws.Range("H135").Value = 284803.25
$ws.Range("J135").Value = 284803.25
$ws.Range("L135").Value = 284803.25
$ws.Range("N135").Value = -294943.25
$ws.Range("H138").Value = 94694.5
$ws.Range("I138").Value = 94390
$ws.Range("J138").Value = 94999
$ws.Range("K138").Value = 94390
$ws.Range("L138").Value = 94999
$ws.Range("M138").Value = -89250
$ws.Range("N138").Value = -105279
